# Regenerate orders with updated distance/size codes.
# The experiment's distance and face-size labels changed:
#   D64 -> D69, D51 -> D55, D80 -> D86 (distances)
#   S30 -> S31 (size)
# These tokens appear embedded inside several text columns (Condition,
# Filename_Left, Filename_Right, Distance, Size), so do a whole-value
# substring replace across the used range for each old->new pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

$replacements = @(
    @("D64", "D69"),
    @("D51", "D55"),
    @("D80", "D86"),
    @("S30", "S31")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $used.Replace($old, $new, -4123, 1, $false, $false, $true)
}
